# Fix the uneven dash spacing around the "Education" and "Work Experience"
# section headings: each heading line currently has 4 dashes before the
# word and too few after it; move one dash from the leading group to the
# trailing group so the line is balanced (same total dash count either
# way - just re-centered).

$d = $word.ActiveDocument
$dashChar = [string][char]0x2500   # U+2500 BOX DRAWINGS LIGHT HORIZONTAL

function Balance-HeadingDashes($paragraphIndex) {
    $para = $d.Paragraphs.Item($paragraphIndex).Range

    # Remove one dash from the run of dashes at the very start of the
    # heading paragraph (e.g. "────" -> "───").
    $lead = $d.Range($para.Start, $para.Start + 1)
    $lead.Delete()

    # Re-fetch the (now one character shorter) paragraph range and find
    # the insertion point for the extra trailing dash: right before the
    # paragraph mark, but before any trailing space that follows the
    # trailing dash run.
    $para2 = $d.Paragraphs.Item($paragraphIndex).Range
    $insertAt = $para2.End - 1
    while ($d.Range($insertAt - 1, $insertAt).Text -ne $dashChar) {
        $insertAt = $insertAt - 1
    }
    $insertPoint = $d.Range($insertAt, $insertAt)
    $insertPoint.InsertBefore($dashChar)
}

# "Education" heading
Balance-HeadingDashes 5

# "Work Experience" heading
Balance-HeadingDashes 10

# The auto "last edit" bookmark should now sit right after the shortened
# leading dash run in the Education heading (this mirrors where Word
# leaves it after this kind of in-place edit).
$d2 = $word.ActiveDocument
$eduPara = $d2.Paragraphs.Item(5).Range
$goBackTarget = $d2.Range($eduPara.Start + 3, $eduPara.Start + 3)
$d2.Bookmarks.Add("_GoBack", $goBackTarget)

Write-Output "Education heading: [$($d2.Paragraphs.Item(5).Range.Text)]"
Write-Output "Work Experience heading: [$($d2.Paragraphs.Item(10).Range.Text)]"
